$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# Correccion de la cantidad de datos vacios: C74 pasa de 5 a 6
$ws.Range("C74").Value = 6

# Copiar el formato de la fila 74 (plantilla) hacia las filas nuevas 75 y 76
$ws.Range("A74:C74").Copy()
$ws.Range("A75:C75").PasteSpecial(-4122)
$ws.Range("A74:C74").Copy()
$ws.Range("A76:C76").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fila 75: nuevo registro "Curso de manejo de datos"
$ws.Range("A75").Value = 44763
$ws.Range("B75").Value = "Curso de manejo de datos"
$ws.Range("C75").Value = 2

# Fila 76: la columna A es texto (no fecha) en este registro
$ws.Range("B74").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A76").Value = "23(07/2022"
$ws.Range("B76").Value = "Reunion con vero organizando el codigo con datos vacios"
$ws.Range("C76").Value = 1.5
$ws.Rows("76:76").RowHeight = 27.6

$ws.Range("C76").Select()
